$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05976570080812333
$ws.Range("H2").Value = -7.0549677135518
$ws.Range("I2").Value = -10.85658129332591
$ws.Range("G3").Value = 0.0679854072668568
$ws.Range("H3").Value = 20.95532421504451
$ws.Range("G4").Value = -0.01987293130034513
$ws.Range("H4").Value = 2.15973419950797
$ws.Range("G5").Value = -0.01942756209528249
$ws.Range("H5").Value = -69.95013228107459
$ws.Range("G6").Value = -0.004185219196956102
$ws.Range("H6").Value = 62.61046469156003
$ws.Range("G7").Value = -0.006273027295040768
$ws.Range("H7").Value = -10.75776276709965
$ws.Range("G8").Value = 0.001637362195990166
$ws.Range("H8").Value = 128.7803706239258
$ws.Range("G9").Value = -0.006602647764756579
$ws.Range("H9").Value = -20.19792445542275
$ws.Range("G10").Value = -0.05961985043809716
$ws.Range("H10").Value = 5.397983971557854
$ws.Range("G11").Value = -0.06011439550248148
$ws.Range("H11").Value = 6.185663592686814
$ws.Range("G12").Value = -0.3924494586741473
$ws.Range("H12").Value = 0.5458498760614592
$ws.Range("G13").Value = -0.3972110673275344
$ws.Range("H13").Value = -1.298550845090402
$ws.Range("G14").Value = -0.029078031878882
$ws.Range("H14").Value = -254.8583935997151
$ws.Range("G15").Value = -0.01641909250093018
$ws.Range("H15").Value = 63.76969905924921
$ws.Range("G16").Value = 0.1392673142470738
$ws.Range("H16").Value = 1.84107595357264
$ws.Range("G17").Value = 0.1471350991494476
$ws.Range("H17").Value = 5.495393664921557
$ws.Range("G18").Value = 0.1162280770196786
$ws.Range("H18").Value = -1.223891168056393
$ws.Range("G19").Value = 0.1301526247866527
$ws.Range("H19").Value = 1.153172299667688
$ws.Range("G20").Value = 0.08466031211345856
$ws.Range("H20").Value = -4.591578798283673
$ws.Range("G21").Value = 0.09460642344713517
$ws.Range("H21").Value = 8.629360093823806
$ws.Range("G22").Value = -0.09604312390261531
$ws.Range("H22").Value = -2.738412178487378
$ws.Range("G23").Value = -0.1075269800455612
$ws.Range("H23").Value = -5.995328340119603
$ws.Range("G24").Value = 0.1637253715976425
$ws.Range("H24").Value = 1.635292539678031
$ws.Range("G25").Value = 0.1731093553342201
$ws.Range("H25").Value = 1.473438603834955
$ws.Range("G26").Value = 0.09614473738344895
$ws.Range("H26").Value = 6.057045546570098
$ws.Range("G27").Value = 0.09553631476917235
$ws.Range("H27").Value = 11.13817355273312
$ws.Range("G28").Value = -0.1347967142038938
$ws.Range("H28").Value = 2.06605139623586
$ws.Range("G29").Value = -0.1321370538260386
$ws.Range("H29").Value = 5.481204934076191
$ws.Range("G30").Value = 0.04097370250001889
$ws.Range("H30").Value = -21.22755705256947
$ws.Range("G31").Value = 0.05076432056760785
$ws.Range("H31").Value = 15.8489552051023
$ws.Range("G32").Value = 0.1135982688903562
$ws.Range("H32").Value = 4.502447289706424
$ws.Range("G33").Value = 0.1151390882572715
$ws.Range("H33").Value = -7.217124212998677
$ws.Range("G34").Value = -0.01539632928691434
$ws.Range("H34").Value = 1.409189988487183
$ws.Range("G35").Value = -0.01137966929216473
$ws.Range("H35").Value = 32.01007824357931
$ws.Range("G36").Value = 0.03612817292176713
$ws.Range("H36").Value = -1.737259422668674
$ws.Range("G37").Value = 0.03104115748035275
$ws.Range("H37").Value = -13.01501074768267
$ws.Range("G38").Value = 0.1071500297611073
$ws.Range("H38").Value = 6.82716668299812
$ws.Range("G39").Value = 0.1036507264547394
$ws.Range("H39").Value = 6.417567274553809
$ws.Range("G40").Value = 0.03079211631859983
$ws.Range("H40").Value = -8.595155876894381
$ws.Range("G41").Value = 0.0334346923000989
$ws.Range("H41").Value = 3.774101130472719
$ws.Range("G42").Value = 0.1187413181626411
$ws.Range("H42").Value = -1.789921523034432
$ws.Range("G43").Value = 0.1273700667740273
$ws.Range("H43").Value = -0.3250608306703475
$ws.Range("G44").Value = 0.03434911578676667
$ws.Range("H44").Value = -13.40553052775557
$ws.Range("G45").Value = 0.03613334281934084
$ws.Range("H45").Value = 15.95190805747183
$ws.Range("G46").Value = 0.05775645585319572
$ws.Range("H46").Value = 2.011370200789755
$ws.Range("G47").Value = 0.06126035046565886
$ws.Range("H47").Value = 4.414787050045385
$ws.Range("G48").Value = 0.05189510672561686
$ws.Range("H48").Value = 5.355110382425165
$ws.Range("G49").Value = 0.04471712144371884
$ws.Range("H49").Value = -1.866727867450029
$ws.Range("G50").Value = 0.02981843696220924
$ws.Range("H50").Value = 12.57982463685091
$ws.Range("G51").Value = 0.02612960981444661
$ws.Range("H51").Value = -6.732822322625085
$ws.Range("G52").Value = -0.08891397169286448
$ws.Range("H52").Value = -2.289745984281968
$ws.Range("G53").Value = -0.0798288179252052
$ws.Range("H53").Value = 0.4829248886763735
$ws.Range("G54").Value = 0.05430404335337361
$ws.Range("H54").Value = 8.550116760178826
$ws.Range("G55").Value = 0.05236120690121006
$ws.Range("H55").Value = -6.971388323905393
$ws.Range("G56").Value = 0.0493017466570904
$ws.Range("H56").Value = -0.2632881671379027
$ws.Range("G57").Value = 0.04703817509318357
$ws.Range("H57").Value = 23.83688530159399
$ws.Range("G58").Value = 0.05698018807881352
$ws.Range("H58").Value = -1.100190842713466
$ws.Range("G59").Value = 0.05547448007740317
$ws.Range("H59").Value = -2.699070136480825
$ws.Range("G60").Value = 0.02899252505027517
$ws.Range("H60").Value = 5.580937941371831
$ws.Range("G61").Value = 0.02512435341370487
$ws.Range("H61").Value = -5.893633224600912
$ws.Range("G62").Value = 0.06381419547026682
$ws.Range("H62").Value = 2.188990881639208
$ws.Range("G63").Value = 0.06309197299284462
$ws.Range("H63").Value = -1.237732188303977
$ws.Range("G64").Value = 0.0314809602244901
$ws.Range("H64").Value = 13.48719900786431
$ws.Range("G65").Value = 0.03331384004207594
$ws.Range("H65").Value = -5.965497622344854
$ws.Range("G66").Value = 0.0740146325386359
$ws.Range("H66").Value = -4.722449019419912
$ws.Range("G67").Value = 0.0849410701845828
$ws.Range("H67").Value = 7.696428186321429
$ws.Range("G68").Value = -0.01933166150326306
$ws.Range("H68").Value = 11.10119258911268
$ws.Range("G69").Value = -0.02300302133597345
$ws.Range("H69").Value = -20.174735671023
$ws.Range("G70").Value = 0.0689355721123537
$ws.Range("H70").Value = -4.237346886254972
$ws.Range("G71").Value = 0.07260196555643064
$ws.Range("H71").Value = -8.5859773115241
$ws.Range("G72").Value = -0.1439372952536587
$ws.Range("H72").Value = 6.316124363532913
$ws.Range("G73").Value = -0.155855605668093
$ws.Range("H73").Value = -1.816129842271878
$ws.Range("G74").Value = 0.1508356937617299
$ws.Range("H74").Value = 0.2768927827326374
$ws.Range("G75").Value = 0.15147854307786
$ws.Range("H75").Value = 0.6854349467786852
$ws.Range("G76").Value = -0.01077430611214037
$ws.Range("H76").Value = -939.5277058118113
$ws.Range("G77").Value = -0.001921176400190697
$ws.Range("H77").Value = 12.98577702939925
$ws.Range("G78").Value = 0.09065826624400664
$ws.Range("H78").Value = 0.7781563383047863
$ws.Range("G79").Value = 0.08880881285479843
$ws.Range("H79").Value = -8.351329099797079
$ws.Range("G80").Value = -0.2237385932061497
$ws.Range("H80").Value = -3.375276840503223
$ws.Range("G81").Value = -0.2219503102922142
$ws.Range("H81").Value = -4.146339377663788
$ws.Range("G82").Value = 0.181620894101806
$ws.Range("H82").Value = 8.354266019297176
$ws.Range("G83").Value = 0.1860061220000803
$ws.Range("H83").Value = 5.667926113939258
$ws.Range("G84").Value = 0.109994194044831
$ws.Range("H84").Value = 3.641296581766572
$ws.Range("G85").Value = 0.1107024372627287
$ws.Range("H85").Value = 5.875765552212708
